# Added a bool which controls heat seeking behavior in rockets,
# and added an augment ("Circuit seeking") that turns it on.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AugmentList")

# --- New augment row (row 11) -------------------------------------------
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Circuit seeking"
$ws.Range("C11").Value = 1

# D11 mirrors the "Code" style used elsewhere in column D (quote-prefixed
# right-aligned text style) - copy the format from D6 (same string value
# "0") then set the value.
$ws.Range("D6").Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4122) | Out-Null
$ws.Range("D11").Value = "0"

$script = "def OnAttached() { " + [char]10 + "AddModifier(`"Rocket`", `"maxClip`", `"Flat`", 1);" + [char]10 + "AddModifier(`"Rocket`",`"heatSeeking`",`"Flat`",3);" + [char]10 + "} "

# F11 / H11 / J11 all hold the same script text, mirror the wrap-text
# formatting used by the other "Code" columns (e.g. F2).
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F11").PasteSpecial(-4122) | Out-Null
$ws.Range("F11").Value = $script

$ws.Range("F2").Copy() | Out-Null
$ws.Range("H11").PasteSpecial(-4122) | Out-Null
$ws.Range("H11").Value = $script

$ws.Range("F2").Copy() | Out-Null
$ws.Range("J11").PasteSpecial(-4122) | Out-Null
$ws.Range("J11").Value = $script

$ws.Rows.Item(11).RowHeight = 60

# --- Sheet view / selection state ---------------------------------------
# Preserve SynergyList's own selection (D3) before switching tabs, since
# selecting a range on another sheet also activates that sheet.
$synergy = $wb.Worksheets.Item("SynergyList")
$synergy.Range("D3").Select() | Out-Null

# AugmentList becomes the active / selected tab (was SynergyList).
$ws.Activate()
$ws.Range("B11").Select() | Out-Null
